$wb = $excel.ActiveWorkbook

$wsPrincipal = $wb.Worksheets.Item("Principal")
$wsDoctores  = $wb.Worksheets.Item("Doctores")
$wsPacientes = $wb.Worksheets.Item("Pacientes")

# --- Sheet: Principal ---
# Fill in the alternate-case row (row 3) that was previously empty
$wsPrincipal.Range("D3").Value = "asanchez"
$wsPrincipal.Range("G3").Value = "Cita Error"
$wsPrincipal.Range("E3").Value = "juribe"

# --- Sheet: Pacientes ---
$wsPacientes.Range("C4").Value = "Uribe"

# --- Sheet: Principal (continued) ---
$wsPrincipal.Range("F3").Value = "20/06/2018"
$wsPrincipal.Range("H3").Value = "Error:"

$wsPrincipal.Range("F3").Select()

# --- Sheet: Doctores ---
$wsDoctores.Range("F4").Value = 10299990

$wsDoctores.Range("F4").Select()

# --- Sheet: Pacientes (continued) ---
$wsPacientes.Range("A4").Value = "juribe"
$wsPacientes.Range("B4").Value = "Juan D."
$wsPacientes.Range("E4").Value = 1010198
$wsPacientes.Range("F4").Value = 10001009

$wsPacientes.Range("A9").Select()

# Restore Principal as the active/selected sheet (tab) without disturbing
# the selections already recorded on the other sheets.
$wsPrincipal.Activate()
